# Update Ligand/Receptor-expressing cell counts and all dependent
# expression/specificity/edge-weight metrics for rows 2-17, per
# "Natmi following Dr Hou advice" (Mdk-Sdc3 LR-pair recompute).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; E = 3; G = 1.324023666666666; H = 3.972071; I = 0.01518042398701374; J = 0.01518042398701374; K = 3; M = 34.36078833333333; N = 103.082365; O = 0.28490270239021; P = 0.28490270239021; Q = 45.49449695865722; R = 409.450472627915; S = 0.004324943817329381; T = 0.00432494381732938 },
    @{ Row = 3; E = 3; G = 1.324023666666666; H = 3.972071; I = 0.01518042398701374; J = 0.01518042398701374; K = 3; M = 21.54461566666667; N = 64.633847; O = 0.17863732245739; P = 0.1786373224573899; Q = 28.52558103190411; R = 256.730229287137; S = 0.002711790294808071; T = 0.00271179029480807 },
    @{ Row = 4; E = 3; G = 1.324023666666666; H = 3.972071; I = 0.01518042398701374; J = 0.01518042398701374; K = 3; M = 60.03138866666666; N = 180.094166; O = 0.4977506538398792; P = 0.4977506538398792; Q = 79.48297933753176; R = 715.346814037786; S = 0.007556065965102674; T = 0.007556065965102674 },
    @{ Row = 5; E = 3; G = 1.324023666666666; H = 3.972071; I = 0.01518042398701374; J = 0.01518042398701374; K = 3; M = 4.668551; N = 14.005653; O = 0.03870932131252084; P = 0.03870932131252084; Q = 6.181272013040332; R = 55.631448117363; S = 0.0005876239097736134; T = 0.0005876239097736134 },
    @{ Row = 6; E = 3; G = 81.17653533333333; H = 243.529606; I = 0.9307191821270077; J = 0.9307191821270075; K = 3; M = 34.36078833333333; N = 103.082365; O = 0.28490270239021; P = 0.28490270239021; Q = 2789.289748222021; R = 25103.60773399819; S = 0.2651644101543906; T = 0.2651644101543905 },
    @{ Row = 7; E = 3; G = 81.17653533333333; H = 243.529606; I = 0.9307191821270077; J = 0.9307191821270075; K = 3; M = 21.54461566666667; N = 64.633847; O = 0.17863732245739; P = 0.1786373224573899; Q = 1748.917254908254; R = 15740.25529417428; S = 0.1662611826549006; T = 0.1662611826549005 },
    @{ Row = 8; E = 3; G = 81.17653533333333; H = 243.529606; I = 0.9307191821270077; J = 0.9307191821270075; K = 3; M = 60.03138866666666; N = 180.094166; O = 0.4977506538398792; P = 0.4977506538398792; Q = 4873.140143208733; R = 43858.2612888786; S = 0.4632660814450357; T = 0.4632660814450356 },
    @{ Row = 9; E = 3; G = 81.17653533333333; H = 243.529606; I = 0.9307191821270077; J = 0.9307191821270075; K = 3; M = 4.668551; N = 14.005653; O = 0.03870932131252084; P = 0.03870932131252084; Q = 378.9767952069686; R = 3410.791156862718; S = 0.03602750787268094; T = 0.03602750787268094 },
    @{ Row = 10; E = 3; G = 1.192675; H = 3.578025; I = 0.0136744626508778; J = 0.0136744626508778; K = 3; M = 34.36078833333333; N = 103.082365; O = 0.28490270239021; P = 0.28490270239021; Q = 40.98125322545834; R = 368.831279029125; S = 0.003895891362969081; T = 0.00389589136296908 },
    @{ Row = 11; E = 3; G = 1.192675; H = 3.578025; I = 0.0136744626508778; J = 0.0136744626508778; K = 3; M = 21.54461566666667; N = 64.633847; O = 0.17863732245739; P = 0.1786373224573899; Q = 25.69572449024167; R = 231.261520412175; S = 0.002442769393996394; T = 0.002442769393996393 },
    @{ Row = 12; E = 3; G = 1.192675; H = 3.578025; I = 0.0136744626508778; J = 0.0136744626508778; K = 3; M = 60.03138866666666; N = 180.094166; O = 0.4977506538398792; P = 0.4977506538398792; Q = 71.59793647801668; R = 644.3814283021501; S = 0.006806472725383435; T = 0.006806472725383434 },
    @{ Row = 13; E = 3; G = 1.192675; H = 3.578025; I = 0.0136744626508778; J = 0.0136744626508778; K = 3; M = 4.668551; N = 14.005653; O = 0.03870932131252084; P = 0.03870932131252084; Q = 5.568064063925001; R = 50.11257657532501; S = 0.0005293291685288943; T = 0.0005293291685288943 },
    @{ Row = 14; E = 3; G = 3.525915333333334; H = 10.577746; I = 0.04042593123510095; J = 0.04042593123510094; K = 3; M = 34.36078833333333; N = 103.082365; O = 0.28490270239021; P = 0.28490270239021; Q = 121.1532304499211; R = 1090.37907404929; S = 0.01151745705552106; T = 0.01151745705552106 },
    @{ Row = 15; E = 3; G = 3.525915333333334; H = 10.577746; I = 0.04042593123510095; J = 0.04042593123510094; K = 3; M = 21.54461566666667; N = 64.633847; O = 0.17863732245739; P = 0.1786373224573899; Q = 75.96449072987357; R = 683.6804165688621; S = 0.007221580113685003; T = 0.007221580113685 },
    @{ Row = 16; E = 3; G = 3.525915333333334; H = 10.577746; I = 0.04042593123510095; J = 0.04042593123510094; K = 3; M = 60.03138866666666; N = 180.094166; O = 0.4977506538398792; P = 0.4977506538398792; Q = 211.6655937810929; R = 1904.990344029836; S = 0.02012203370435749; T = 0.02012203370435749 },
    @{ Row = 17; E = 3; G = 3.525915333333334; H = 10.577746; I = 0.04042593123510095; J = 0.04042593123510094; K = 3; M = 4.668551; N = 14.005653; O = 0.03870932131252084; P = 0.03870932131252084; Q = 16.46091555534867; R = 148.148239998138; S = 0.001564860361537395; T = 0.001564860361537395 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 7).Value = $u.G
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
    $ws.Cells.Item($u.Row, 11).Value = $u.K
    $ws.Cells.Item($u.Row, 13).Value = $u.M
    $ws.Cells.Item($u.Row, 14).Value = $u.N
    $ws.Cells.Item($u.Row, 15).Value = $u.O
    $ws.Cells.Item($u.Row, 16).Value = $u.P
    $ws.Cells.Item($u.Row, 17).Value = $u.Q
    $ws.Cells.Item($u.Row, 18).Value = $u.R
    $ws.Cells.Item($u.Row, 19).Value = $u.S
    $ws.Cells.Item($u.Row, 20).Value = $u.T
}
